$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "4、8次情侣畅玩卡" deposit-card row (old row 22), which duplicated
# the codeAndDeposit/booking figures that now live on their own rows. All
# subsequent rows shift up by one.
$ws.Rows(22).Delete()

# Update the active selection to match the post-edit cursor position.
$ws.Range("B6").Select()
